$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AD, AE, AF) in row 1,
# matching the header style used by the existing headers (copy format from AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122, $null, $false, $false) | Out-Null
$excel.CutCopyMode = 0

# Data rows 2-49: team record values (Wins=78, Losses=84, Ties=0)
$data = $ws.Range("AD2:AF49")
$values = New-Object 'object[,]' 48,3
for ($i = 0; $i -lt 48; $i++) {
    $values[$i, 0] = 78
    $values[$i, 1] = 84
    $values[$i, 2] = 0
}
$data.Value = $values
